# "Doing Updates for Financials" — add the newest fiscal-year column (FY
# ending 2018-12-31, serial 43465) to the PLOW yearly financial-statement
# sheet. The new column is inserted immediately before the existing column
# D, pushing the prior D:K data right into E:L, matching how the source
# spreadsheet is refreshed each year with one more period of history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D; Excel shifts D:K -> E:L.
$ws.Columns("D").Insert()

# The freshly inserted column has no formatting of its own yet - copy it
# from the (now-shifted) column E so the new D cells pick up the same
# date / number styles as the rest of each row.
$fmtSrc = $ws.Range("E7:E102")
$fmtDst = $ws.Range("D7:D102")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# New period-ending header (2018-12-31) for the three statement blocks.
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# --- Income Statement --------------------------------------------------
$ws.Range("D8").Value = 524100    # Total Revenue
$ws.Range("D9").Value = 369200    # Cost of Revenue
$ws.Range("D10").Value = 154900   # Gross Profit
$ws.Range("D12").Value = "NA"     # Research Development
$ws.Range("D13").Value = 0        # NA
$ws.Range("D14").Value = 0        # Selling General and Administrative
$ws.Range("D15").Value = 11500    # Non Recurring
$ws.Range("D17").Value = 450600   # Total Operating Expenses
$ws.Range("D18").Value = 73500    # Operating Income or Loss
$ws.Range("D20").Value = -17700   # Total Other Income/Expenses Net
$ws.Range("D21").Value = 74800    # Earnings Before Interest And Taxes
$ws.Range("D22").Value = 0        # Interest Expense
$ws.Range("D23").Value = 55800    # Income Before Tax
$ws.Range("D24").Value = 11900    # Income Tax Expense
$ws.Range("D25").Value = 0        # Minority Interest
$ws.Range("D26").Value = 43900    # Income After Tax
$ws.Range("D27").Value = 43300    # Net Income From Continuing Ops
$ws.Range("D28").Value = 0        # Non-recurring Events
$ws.Range("D29").Value = 0        # Discontinued Operations
$ws.Range("D30").Value = 0        # Extraordinary Items
$ws.Range("D31").Value = 0        # Effect Of Accounting Changes
$ws.Range("D32").Value = 17700    # Other Items
$ws.Range("D33").Value = 43300    # Net Income
$ws.Range("D34").Value = 0        # Preferred Stock And Other Adjustments
$ws.Range("D35").Value = 43300    # Net Income Applicable To Common Shares

# --- Balance Sheet -------------------------------------------------------
$ws.Range("D41").Value = 27800    # Cash And Cash Equivalents
$ws.Range("D42").Value = 0        # Short Term Investments
$ws.Range("D43").Value = 81500    # Net Receivables
$ws.Range("D44").Value = 86200    # Inventory
$ws.Range("D45").Value = 3600     # Other Current Assets
$ws.Range("D46").Value = 199100   # Total Current Assets
$ws.Range("D47").Value = 0        # Long Term Investments
$ws.Range("D48").Value = 55200    # Property Plant and Equipment
$ws.Range("D49").Value = 415700   # Goodwill
$ws.Range("D50").Value = 0        # Intangible Assets
$ws.Range("D51").Value = 0        # Accumulated Amortization
$ws.Range("D52").Value = 6200     # Other Assets
$ws.Range("D53").Value = 0        # Deferred Long Term Asset Charges
$ws.Range("D54").Value = 676200   # Total Assets
$ws.Range("D57").Value = 18700    # Accounts Payable
$ws.Range("D58").Value = 32700    # Short/Current Long Term Debt
$ws.Range("D59").Value = 27600    # Other Current Liabilities
$ws.Range("D60").Value = 79100    # Total Current Liabilities
$ws.Range("D61").Value = 242900   # Long Term Debt
$ws.Range("D62").Value = 71400    # Other Liabilities
$ws.Range("D63").Value = 0        # Deferred Long Term Liability Charges
$ws.Range("D64").Value = 0        # Negative Goodwill
$ws.Range("D65").Value = 0        # Total Liabilities (NA row placeholder)
$ws.Range("D66").Value = 393400   # Total Liabilities
$ws.Range("D68").Value = 0        # Misc Stocks Options Warrants
$ws.Range("D69").Value = 0        # Redeemable Preferred Stock
$ws.Range("D70").Value = 0        # Preferred Stock
$ws.Range("D71").Value = 0        # Common Stock
$ws.Range("D72").Value = 136800   # Retained Earnings
$ws.Range("D73").Value = 0        # Treasury Stock
$ws.Range("D74").Value = 0        # Capital Surplus
$ws.Range("D75").Value = 0        # Other Stockholder Equity
$ws.Range("D76").Value = 282800   # Total Stockholder Equity
$ws.Range("D77").Value = 0        # Net Tangible Assets

# --- Cash Flow Statement ---------------------------------------------------
$ws.Range("D81").Value = 43300    # Net Income
$ws.Range("D83").Value = 19100    # Depreciation
$ws.Range("D84").Value = 0        # Adjustments To Net Income
$ws.Range("D85").Value = 0        # Changes In Accounts Receivables
$ws.Range("D86").Value = 0        # Changes In Liabilities
$ws.Range("D87").Value = 0        # Changes In Inventories
$ws.Range("D88").Value = 0        # Changes In Other Operating Activities
$ws.Range("D89").Value = 58200    # Total Cash Flow From Operating Activities
$ws.Range("D91").Value = -9700    # Capital Expenditures
$ws.Range("D92").Value = 0        # Investments
$ws.Range("D93").Value = 0        # Other Cashflows from Investing Activities
$ws.Range("D94").Value = -9700    # Total Cash Flows From Investing Activities
$ws.Range("D96").Value = -24400   # Dividends Paid
$ws.Range("D97").Value = 0        # Sale Purchase of Stock
$ws.Range("D98").Value = 0        # Net Borrowings
$ws.Range("D99").Value = 0        # Other Cash Flows from Financing Activities
$ws.Range("D100").Value = -57500  # Total Cash Flows From Financing Activities
$ws.Range("D101").Value = 0       # Effect Of Exchange Rate Changes
$ws.Range("D102").Value = -9100   # Change In Cash and Cash Equivalents

Write-Output "Inserted new FY2018 column D and populated financial data"
